# merged with latest from UW team's update
# The form_id setting in the "settings" sheet changes from "refrigerators"
# to "refrigerators_update", while table_id stays "refrigerators".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Range("B2").Value = "refrigerators_update"
